# CDS Study filter fixes
#
# The "ParticipantsTab" row's Cypher query (cell B2 on the "startup" sheet)
# is replaced with an updated/fixed query (adds genomic_info traversal,
# reversed participant/study relationship direction, sorted sample list,
# reformatted ORDER BY / LIMIT). The row grows taller to fit the longer
# text, and the active selection moves from A2 to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["Clonal evolution during metastatic spread in high-rish neuroblastoma"]
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id
LIMIT 100
'@

$ws.Range("B2").Value = $newQuery

# The new query text wraps onto more lines, so the row grows to fit it.
$ws.Rows.Item(2).RowHeight = 283.5

# The author's selection ended up on the query cell they just edited.
[void]$ws.Range("B2").Select()
